# Updates the cryptocurrency price/volume table with the latest scraped values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value looks like a plain number need to be forced to Text
# format first, otherwise Excel will auto-convert them to numeric values and
# lose formatting (trailing zeros) or introduce floating point noise.
$textFormatCells = @(
    'D5', 'D6', 'D7', 'D8', 'D9', 'D10', 'D12', 'D13',
    'D14', 'D16', 'D19', 'D21', 'D22', 'D23', 'D25', 'D26',
    'D27', 'D28', 'D29', 'D30', 'D31', 'D33', 'D34', 'D35',
    'D37', 'D38', 'D39', 'D41', 'D42', 'D43', 'D45', 'D47',
    'D48', 'D49', 'D50', 'D51'
)
foreach ($cellRef in $textFormatCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

# Apply the updated values cell-by-cell, row by row.
# Row 2
$ws.Range('D2').Value = '36.968.25'
$ws.Range('E2').Value = '  +4.49%  '
# Row 3
$ws.Range('D3').Value = '1.915.37'
$ws.Range('E3').Value = '  +1.55%  '
# Row 4
$ws.Range('E4').Value = '  +0.03%  '
# Row 5
$ws.Range('D5').Value = '248.77'
$ws.Range('E5').Value = '  +1.04%  '
# Row 6
$ws.Range('D6').Value = '0.689'
$ws.Range('E6').Value = '  -0.72%  '
# Row 7
$ws.Range('D7').Value = '0.999'
$ws.Range('E7').Value = '  -0.04%  '
# Row 8
$ws.Range('D8').Value = '47.94'
$ws.Range('E8').Value = '  +10.62%  '
# Row 9
$ws.Range('D9').Value = '0.374'
$ws.Range('E9').Value = '  +5.52%  '
# Row 10
$ws.Range('D10').Value = '58.08'
$ws.Range('E10').Value = '  +6.68%  '
# Row 11
$ws.Range('E11').Value = '  +2.03%  '
# Row 12
$ws.Range('D12').Value = '0.100'
$ws.Range('E12').Value = '  +2.19%  '
# Row 13
$ws.Range('D13').Value = '15.23'
$ws.Range('E13').Value = '  +11.50%  '
# Row 14
$ws.Range('D14').Value = '0.822'
$ws.Range('E14').Value = '  +7.31%  '
# Row 15
$ws.Range('D15').Value = '2.192.32'
$ws.Range('E15').Value = '  +1.55%  '
# Row 16
$ws.Range('D16').Value = '5.12'
$ws.Range('E16').Value = '  +1.74%  '
# Row 17
$ws.Range('D17').Value = '1.913.03'
$ws.Range('E17').Value = '  +1.50%  '
# Row 18
$ws.Range('D18').Value = '37.051.41'
$ws.Range('E18').Value = '  +4.88%  '
# Row 19
$ws.Range('D19').Value = '74.55'
$ws.Range('E19').Value = '  +1.56%  '
# Row 20
$ws.Range('D20').Value = '0.0₃0855'
$ws.Range('E20').Value = '  +3.39%  '
# Row 21
$ws.Range('D21').Value = '13.69'
$ws.Range('E21').Value = '  +7.09%  '
# Row 22
$ws.Range('D22').Value = '251.24'
$ws.Range('E22').Value = '  +2.84%  '
# Row 23
$ws.Range('D23').Value = '5.15'
$ws.Range('E23').Value = '  -0.55%  '
# Row 24
$ws.Range('E24').Value = '  +0.14%  '
# Row 25
$ws.Range('D25').Value = '2.43'
$ws.Range('E25').Value = '  -7.80%  '
# Row 26
$ws.Range('D26').Value = '2.18'
$ws.Range('E26').Value = '  +2.23%  '
# Row 27
$ws.Range('D27').Value = '167.40'
$ws.Range('E27').Value = '  +1.34%  '
# Row 28
$ws.Range('D28').Value = '8.83'
$ws.Range('E28').Value = '  +2.25%  '
# Row 29
$ws.Range('D29').Value = '18.71'
$ws.Range('E29').Value = '  +2.42%  '
# Row 30
$ws.Range('D30').Value = '0.129'
$ws.Range('E30').Value = '  +0.63%  '
# Row 31
$ws.Range('D31').Value = '4.59'
$ws.Range('E31').Value = '  +7.06%  '
# Row 32
$ws.Range('E32').Value = '  +2.60%  '
# Row 33
$ws.Range('D33').Value = '4.31'
$ws.Range('E33').Value = '  +2.98%  '
# Row 34
$ws.Range('D34').Value = '0.0900'
$ws.Range('E34').Value = '  +23.42%  '
# Row 35
$ws.Range('D35').Value = '1.90'
$ws.Range('E35').Value = '  +1.87%  '
# Row 36
$ws.Range('E36').Value = '  +0.04%  '
# Row 37
$ws.Range('D37').Value = '19.40'
$ws.Range('E37').Value = '  +57.95%  '
# Row 38
$ws.Range('D38').Value = '1.48'
$ws.Range('E38').Value = '  +2.03%  '
# Row 39
$ws.Range('D39').Value = '0.884'
$ws.Range('E39').Value = '  +3.32%  '
# Row 40
$ws.Range('E40').Value = '  +1.88%  '
# Row 41
$ws.Range('D41').Value = '104.60'
$ws.Range('E41').Value = '  +7.50%  '
# Row 42
$ws.Range('D42').Value = '0.0227'
$ws.Range('E42').Value = '  +3.92%  '
# Row 43
$ws.Range('D43').Value = '17.61'
$ws.Range('E43').Value = '  +2.08%  '
# Row 44
$ws.Range('E44').Value = '  +19.97%  '
# Row 45
$ws.Range('D45').Value = '1.09'
$ws.Range('E45').Value = '  +2.17%  '
# Row 46
$ws.Range('D46').Value = '1.353.36'
$ws.Range('E46').Value = '  +3.41%  '
# Row 47
$ws.Range('D47').Value = '2.39'
$ws.Range('E47').Value = '  +0.04%  '
# Row 48
$ws.Range('D48').Value = '0.0827'
$ws.Range('E48').Value = '  +1.93%  '
# Row 49
$ws.Range('D49').Value = '2.82'
$ws.Range('E49').Value = '  +3.01%  '
# Row 50
$ws.Range('D50').Value = '6.41'
$ws.Range('E50').Value = '  +1.60%  '
# Row 51
$ws.Range('B51').Value = 'THORChain'
$ws.Range('C51').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D51').Value = '3.76'
$ws.Range('E51').Value = '  +12.68%  '

Write-Host "Updated cryptos list."
